$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A28").Value = "verwalten"

$ws.Range("E9").Select()
